# Generate Report for Handoff
#
# Refreshes the localization-status report: the handoff xliff for
# 2cdba8a9-dd3c-4817-8243-f9701655a2be.md was (re)generated, so its
# priority flips from "low" to "ht" on both language tabs, the zh-cn
# "Latest Handoff Datetime" advances to the new generation time, and the
# Overview tab's "Latest HO Xliff Generate Date" advances as well.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("H4").Value = "2016-08-30 04:30:56"
$wsZh.Range("H5").Value = "2016-08-30 04:30:56"
$wsZh.Range("H6").Value = "2016-08-30 04:30:56"
$wsZh.Range("H7").Value = "2016-08-30 04:30:56"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("E7").Value = "ht"

$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("G4").Value = "2016-08-30 04:31:03"
$wsOv.Range("G5").Value = "2016-08-30 04:31:03"
$wsOv.Range("G6").Value = "2016-08-30 04:31:03"
$wsOv.Range("G7").Value = "2016-08-30 04:31:03"
